$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.870.36"
$ws.Range("E2").Value = "  +4.46%  "
$ws.Range("D3").Value = "3.142.39"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.67"
$ws.Range("E5").Value = "  +2.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.80"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.136.93"
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +19.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  +5.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  +8.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.88"
$ws.Range("E14").Value = "  +4.05%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "3.667.27"
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").Value = "63.814.51"
$ws.Range("E17").Value = "  +4.34%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.17"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.142.91"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.95"
$ws.Range("E20").Value = "  +5.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  +4.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.37"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.45"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.63"
$ws.Range("E27").Value = "  +7.50%  "
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.85"
$ws.Range("E31").Value = "  +6.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.04"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("D34").Value = "0.0₃0869"
$ws.Range("E34").Value = "  +8.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").Value = "  +12.00%  "
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("E37").Value = "  +16.40%  "
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.97"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "449.12"
$ws.Range("E40").Value = "  +9.78%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0377"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("D43").Value = "2.920.69"
$ws.Range("E43").Value = "  +4.87%  "
$ws.Range("E44").Value = "  +7.77%  "
$ws.Range("E45").Value = "  +4.68%  "
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.22"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.84"
$ws.Range("E50").Value = "  +4.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.92"
$ws.Range("E51").Value = "  -9.23%  "
